$wb = $excel.ActiveWorkbook

# --- Update status + timestamps for the handoff-ready report ---

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column(s): "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date (Overview) / Latest Handback DateTime (de-de)
$overview.Range("G2").Value = "2016-08-26 02:58:51"
$dede.Range("H2").Value = "2016-08-26 02:58:51"

# Latest Handoff Datetime (zh-cn)
$zhcn.Range("H2").Value = "2016-08-26 02:58:46"

# --- Column widths grew to fit the new, longer "Ready for handoff" text ---
# (ColumnWidth is specified in "characters"; Excel stores the sheet's <col>
# width in a slightly different, pixel-quantized scale, so the value below
# is the character-width input that lands closest to the saved width.)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
